$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    "B2"=15.38346847130473; "C2"=7.280544208531666; "D2"=6.055608943476014; "E2"=11.91086431536556; "G2"=66.50200929565972; "H2"=23.26262981018506; "K2"=11.74920845148932; "L2"=10.12821394648012; "M2"=15.73052549114441
    "B3"=15.30741413924237; "C3"=7.217209755070515; "D3"=5.949011235723607; "E3"=11.91434090046809; "G3"=65.73936761230179; "H3"=23.16442790394495; "K3"=11.71458864186634; "L3"=10.13973818284178; "M3"=15.74535004385269
    "B4"=15.26571530231473; "C4"=7.176646755053707; "D4"=5.884459735914807; "E4"=11.91742384518966; "G4"=65.27427657860711; "H4"=23.10615278455427; "K4"=11.697429111409; "L4"=10.14822343438494; "M4"=15.75765658065341
    "B5"=15.24999626968472; "C5"=7.159696210582239; "D5"=5.858418593625465; "E5"=11.91891893549789; "G5"=65.08569954138585; "H5"=23.0829230476965; "K5"=11.69147340558386; "L5"=10.1520358042833; "M5"=15.76347740914248
    "B6"=15.24746350080014; "C6"=7.156856186138648; "D6"=5.854111535101054; "E6"=11.91918162331141; "G6"=65.05444835294465; "H6"=23.0790973122014; "K6"=11.69054726128167; "L6"=10.15269026412983; "M6"=15.76449262011929
    "B7"=15.26549813294288; "C7"=7.176419857143007; "D7"=5.884107416709463; "E7"=11.91744304135824; "G7"=65.27172931044473; "H7"=23.10583739134967; "K7"=11.69734458439877; "L7"=10.14827341354336; "M7"=15.75773181975839
    "B8"=15.35621801892016; "C8"=7.259054491909724; "D8"=6.018690002533626; "E8"=11.91186641047746; "G8"=66.23848425260593; "H8"=23.22835620715694; "K8"=11.73642560493706; "L8"=10.13189511640592; "M8"=15.73497209314916
    "B9"=15.57295472600438; "C9"=7.407791231284778; "D9"=6.288085734798439; "E9"=11.90844162672656; "G9"=68.15266169559034; "H9"=23.48420979387858; "K9"=11.84521292168214; "L9"=10.11095094396574; "M9"=15.71575129746281
    "B10"=15.75462195375718; "C10"=7.508941384810345; "D10"=6.487156493430517; "E10"=11.91048630107894; "G10"=69.56104018032231; "H10"=23.68106864910666; "K10"=11.94416281100519; "L10"=10.10236166170082; "M10"=15.71709022764091
    "B11"=15.84183585341287; "C11"=7.553191085169376; "D11"=6.577519345354621; "E11"=11.91240243474114; "G11"=70.20024663608449; "H11"=23.77241003029418; "K11"=11.99315841865515; "L11"=10.09992647656447; "M11"=15.72104419807414
    "B12"=15.87549138891311; "C12"=7.569693419000085; "D12"=6.611673618744511; "E12"=11.91326935124676; "G12"=70.44192457230845; "H12"=23.80724284476705; "K12"=12.0122700233172; "L12"=10.09921557206447; "M12"=15.72302088022226
    "B13"=15.86821552838459; "C13"=7.566150651763983; "D13"=6.604321296919391; "E13"=11.91307636792187; "G13"=70.38989419475246; "H13"=23.79973035047864; "K13"=12.00812942444405; "L13"=10.09935929050122; "M13"=15.72257387348546
    "B14"=15.84459226293783; "C14"=7.554553836127996; "D14"=6.580330725272685; "E14"=11.91247092705801; "G14"=70.22013808410334; "H14"=23.77527093926995; "K14"=11.99471964620491; "L14"=10.09986375907033; "M14"=15.72119722321516
    "B15"=15.83020346020969; "C15"=7.547417353496601; "D15"=6.56562640551339; "E15"=11.91211846638037; "G15"=70.11610381946315; "H15"=23.76032018795658; "K15"=11.98657798688024; "L15"=10.10020025706616; "M15"=15.72041636474284
    "B16"=15.74901196821596; "C16"=7.50601403065451; "D16"=6.481244017254481; "E16"=11.91038087894843; "G16"=69.51922283581629; "H16"=23.67513413794673; "K16"=11.94103974705748; "L16"=10.10255039328399; "M16"=15.71689900490278
    "B17"=15.7003552620157; "C17"=7.480162323503444; "D17"=6.429401954042138; "E17"=11.90956714248788; "G17"=69.15256118492175; "H17"=23.623324088524; "K17"=11.91411351405774; "L17"=10.10436887693656; "M17"=15.71559688761239
    "B18"=15.67280136965803; "C18"=7.465127089115206; "D18"=6.399567442662906; "E18"=11.90919194010825; "G18"=68.94153822427262; "H18"=23.59369357257443; "K18"=11.89900209200522; "L18"=10.10555341421942; "M18"=15.71516289494767
    "B19"=15.66354715016098; "C19"=7.460007916837243; "D19"=6.389464309489172; "E19"=11.90908086158574; "G19"=68.87007238011203; "H19"=23.58369066196883; "K19"=11.8939506031285; "L19"=10.10597829420162; "M19"=15.71507008907721
    "B20"=15.70549032351453; "C20"=7.482931466933536; "D20"=6.434922586209773; "E20"=11.90964416108618; "G20"=69.19160747842098; "H20"=23.62882193597296; "K20"=11.91694106295204; "L20"=10.10416095536097; "M20"=15.71570291464226
    "B21"=15.85151413306986; "C21"=7.557967000279277; "D21"=6.587379359219573; "E21"=11.91264492828049; "G21"=70.27001105791145; "H21"=23.78244875378316; "K21"=11.99864340092578; "L21"=10.09970985506929; "M21"=15.72158858218661
    "B22"=15.95060462933765; "C22"=7.60552600853535; "D22"=6.686629000117823; "E22"=11.91542957222363; "G22"=70.97254642457639; "H22"=23.88426735031659; "K22"=12.05528565487359; "L22"=10.09803186465519; "M22"=15.72822856088128
    "B23"=15.89739326288698; "C23"=7.580278507382214; "D23"=6.633704761004057; "E23"=11.91386817145921; "G23"=70.59785064642269; "H23"=23.82979997988434; "K23"=12.024762832281; "L23"=10.09881495995728; "M23"=15.72442969092654
    "B24"=15.70316745342253; "C24"=7.48168007524673; "D24"=6.432426798547375; "E24"=11.90960905243687; "G24"=69.17395533521506; "H24"=23.62633587285217; "K24"=11.91566157813934; "L24"=10.10425452350153; "M24"=15.71565399972816
    "B25"=15.51028533392204; "C25"=7.368980385428903; "D25"=6.214849854938401; "E25"=11.90856574921685; "G25"=67.63380494419594; "H25"=23.41338724845117; "K25"=11.81239533352901; "L25"=10.11542183222757; "M25"=15.71823265917943
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
